$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for columns I and J, matching the formatting of the existing header cells (e.g. H1)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:I96 and J2:J96
$iValues = @(10,8,8,7,6,6,8,8,8,8,8,8,8,10,8,8,8,8,8,8,8,8,8,8,8,7,7,8,7,8,8,7,8,9,7,8,8,7,8,8,8,8,8,8,8,7,8,8,8,8,8,9,8,8,8,7,8,8,8,8,8,7,7,8,7,7,8,8,8,7,10,8,8,6,8,8,8,9,7,7,7,6,7,7,7,7,7,7,5,6,8,7,4,3,2)
$jValues = @(10,8,8,7,7,7,8,8,8,8,8,8,8,10,8,8,8,8,8,8,8,8,8,8,8,7,7,8,8,8,8,7,8,9,7,8,8,7,8,8,8,8,8,8,8,7,8,8,8,8,8,9,8,8,8,7,8,8,8,8,8,8,7,8,7,7,8,9,8,7,10,8,8,6,8,8,8,9,8,7,7,6,7,7,7,7,7,8,5,6,8,7,5,4,2)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
